$d = $word.ActiveDocument

# Remove the two paragraphs that hold the TBD tags:
#   "[PUMP:TBD:1]"  and  "BOLUS:SRS:2" (ListBullet style)
# These are paragraphs 3 and 4 in the document body (after the Title
# paragraph and the blank paragraph that follow it).
$startPara = $d.Paragraphs.Item(3)
$endPara   = $d.Paragraphs.Item(4)

$range = $d.Range($startPara.Range.Start, $endPara.Range.End)
$range.Delete()
